$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 208, pushing the existing rows 208:295 down to 209:296
$ws.Rows(208).Insert()

# Populate the newly inserted row 208 with its data
$ws.Range("A208").Value2() = 9
$ws.Range("B208").Value2() = 'Vega Central Mapocho de Santiago'
$ws.Range("C208").Value2() = 'Metropolitana'
$ws.Range("D208").Value2() = 44784
$ws.Range("E208").Value2() = 13
$ws.Range("F208").Value2() = 100112001
$ws.Range("G208").Value2() = 'Berenjena'
$ws.Range("H208").Value2() = 'Sin especificar'
$ws.Range("I208").Value2() = 'Primera'
$ws.Range("J208").Value2() = 160
$ws.Range("K208").Value2() = 10000
$ws.Range("L208").Value2() = 10000
$ws.Range("M208").Value2() = 10000
$ws.Range("N208").Value2() = '$/caja 50 unidades'
$ws.Range("O208").Value2() = 'Región de Arica y Parinacota'
$ws.Range("P208").Value2() = 200
$ws.Range("Q208").Value2() = 50
$ws.Range("R208").Value2() = 'Hortaliza'
